$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-scrape dropped a duplicate "10X Bingo Multiplier" ($3.00 Games, game #972) row;
# delete it and let every following row shift up by one.
$ws.Rows.Item(18).Delete()

# Force column F (LAST SCRAPE DATE) to store literal text so Excel does not
# reinterpret "yyyy-mm-dd" strings as date serials.
$ws.Columns.Item(6).NumberFormat = "@"

# Refresh game name / number / remaining-prizes / scrape-date per the latest scrape
$ws.Cells.Item(2,3).Value = "Celebrate!"
$ws.Cells.Item(2,4).Value = 885
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = "2019-03-12"
$ws.Cells.Item(3,3).Value = "7-11-21®"
$ws.Cells.Item(3,4).Value = 953
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = "2019-03-12"
$ws.Cells.Item(4,3).Value = "Blackjack"
$ws.Cells.Item(4,4).Value = 1003
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = "2019-03-12"
$ws.Cells.Item(5,3).Value = "Pocket Change"
$ws.Cells.Item(5,4).Value = 996
$ws.Cells.Item(5,5).Value = 19
$ws.Cells.Item(5,6).Value = "2019-03-12"
$ws.Cells.Item(6,3).Value = "Quick 7s"
$ws.Cells.Item(6,4).Value = 982
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = "2019-03-12"
$ws.Cells.Item(7,3).Value = "Pocket Change"
$ws.Cells.Item(7,4).Value = 970
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = "2019-03-12"
$ws.Cells.Item(8,3).Value = "Merry Money!"
$ws.Cells.Item(8,4).Value = 1011
$ws.Cells.Item(8,5).Value = 56
$ws.Cells.Item(8,6).Value = "2019-03-12"
$ws.Cells.Item(9,3).Value = "Super Triple 7s"
$ws.Cells.Item(9,4).Value = 1014
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = "2019-03-12"
$ws.Cells.Item(10,3).Value = "I Love Cash"
$ws.Cells.Item(10,4).Value = 1022
$ws.Cells.Item(10,5).Value = 5
$ws.Cells.Item(10,6).Value = "2019-03-12"
$ws.Cells.Item(11,3).Value = "Truck$ & Buck$™"
$ws.Cells.Item(11,4).Value = 1015
$ws.Cells.Item(11,5).Value = 6
$ws.Cells.Item(11,6).Value = "2019-03-12"
$ws.Cells.Item(12,3).Value = "Ben There Win That"
$ws.Cells.Item(12,4).Value = 992
$ws.Cells.Item(12,5).Value = 2347
$ws.Cells.Item(12,6).Value = "2019-03-12"
$ws.Cells.Item(13,3).Value = "Bonus Crossword"
$ws.Cells.Item(13,4).Value = 986
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = "2019-03-12"
$ws.Cells.Item(14,3).Value = "10X Bingo Multiplier"
$ws.Cells.Item(14,4).Value = 1030
$ws.Cells.Item(14,5).Value = 4
$ws.Cells.Item(14,6).Value = "2019-03-12"
$ws.Cells.Item(15,3).Value = "The Voice"
$ws.Cells.Item(15,4).Value = 993
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = "2019-03-12"
$ws.Cells.Item(16,3).Value = "Lucky Symbols LOTERIA"
$ws.Cells.Item(16,4).Value = 938
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = "2019-03-12"
$ws.Cells.Item(17,3).Value = "Lucky Ladybug Crossword"
$ws.Cells.Item(17,4).Value = 925
$ws.Cells.Item(17,5).Value = 18
$ws.Cells.Item(17,6).Value = "2019-03-12"
$ws.Cells.Item(18,3).Value = "Red White & Blue"
$ws.Cells.Item(18,4).Value = 985
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = "2019-02-27"
$ws.Cells.Item(19,3).Value = "Triple Tripler"
$ws.Cells.Item(19,4).Value = 971
$ws.Cells.Item(19,5).Value = 2
$ws.Cells.Item(19,6).Value = "2019-03-12"
$ws.Cells.Item(20,3).Value = "Wild Cherry Crossword Doubler"
$ws.Cells.Item(20,4).Value = 998
$ws.Cells.Item(20,5).Value = 1
$ws.Cells.Item(20,6).Value = "2019-03-12"
$ws.Cells.Item(21,3).Value = "N"
$ws.Cells.Item(21,4).Value = 931
$ws.Cells.Item(21,5).Value = 2
$ws.Cells.Item(21,6).Value = "2019-03-12"
$ws.Cells.Item(22,3).Value = "Caesars®"
$ws.Cells.Item(22,4).Value = 1017
$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(22,6).Value = "2019-03-12"
$ws.Cells.Item(23,3).Value = "Power 5"
$ws.Cells.Item(23,4).Value = 1009
$ws.Cells.Item(23,5).Value = 2
$ws.Cells.Item(23,6).Value = "2019-03-12"
$ws.Cells.Item(24,3).Value = "$500 Frenzy"
$ws.Cells.Item(24,4).Value = 979
$ws.Cells.Item(24,5).Value = 67
$ws.Cells.Item(24,6).Value = "2019-03-12"
$ws.Cells.Item(25,3).Value = "Lucky Break"
$ws.Cells.Item(25,4).Value = 1012
$ws.Cells.Item(25,5).Value = 14
$ws.Cells.Item(25,6).Value = "2019-03-12"
$ws.Cells.Item(26,3).Value = "Black Diamond Crossword"
$ws.Cells.Item(26,4).Value = 990
$ws.Cells.Item(26,5).Value = 1
$ws.Cells.Item(26,6).Value = "2019-03-05"
$ws.Cells.Item(27,3).Value = "The Big $10 Super Ticket ®"
$ws.Cells.Item(27,4).Value = 1016
$ws.Cells.Item(27,5).Value = 3
$ws.Cells.Item(27,6).Value = "2019-03-12"
$ws.Cells.Item(28,3).Value = "$2,000,000 Explosion!"
$ws.Cells.Item(28,4).Value = 995
$ws.Cells.Item(28,5).Value = 16
$ws.Cells.Item(28,6).Value = "2019-03-12"
$ws.Cells.Item(29,3).Value = "Double Diamonds"
$ws.Cells.Item(29,4).Value = 914
$ws.Cells.Item(29,5).Value = 1
$ws.Cells.Item(29,6).Value = "2019-03-12"
$ws.Cells.Item(30,3).Value = "$100,000 Cash Crossword"
$ws.Cells.Item(30,4).Value = 1021
$ws.Cells.Item(30,5).Value = 3
$ws.Cells.Item(30,6).Value = "2019-03-12"
$ws.Cells.Item(31,3).Value = "Nebraska Cash Blowout"
$ws.Cells.Item(31,4).Value = 878
$ws.Cells.Item(31,5).Value = 2341
$ws.Cells.Item(31,6).Value = "2019-03-12"
$ws.Cells.Item(32,3).Value = "Winter Takes It All"
$ws.Cells.Item(32,4).Value = 967
$ws.Cells.Item(32,5).Value = 2
$ws.Cells.Item(32,6).Value = "2019-02-25"
$ws.Cells.Item(33,3).Value = "$1,500 Frenzy"
$ws.Cells.Item(33,4).Value = 980
$ws.Cells.Item(33,5).Value = 28
$ws.Cells.Item(33,6).Value = "2019-02-13"
$ws.Cells.Item(34,3).Value = "Multitude Of Money"
$ws.Cells.Item(34,4).Value = 930
$ws.Cells.Item(34,5).Value = 3
$ws.Cells.Item(34,6).Value = "2019-03-12"
$ws.Cells.Item(35,3).Value = "Ultimate Bonus Crossword"
$ws.Cells.Item(35,4).Value = 999
$ws.Cells.Item(35,5).Value = 2
$ws.Cells.Item(35,6).Value = "2019-03-12"
$ws.Cells.Item(36,3).Value = "Multiplier Spectacular"
$ws.Cells.Item(36,4).Value = 1005
$ws.Cells.Item(36,5).Value = 5
$ws.Cells.Item(36,6).Value = "2019-03-12"
$ws.Cells.Item(37,3).Value = "Flawless Fortune"
$ws.Cells.Item(37,4).Value = 890
$ws.Cells.Item(37,5).Value = 1
$ws.Cells.Item(37,6).Value = "2019-03-12"

# Drop the temporary text format now that the literal date strings are committed,
# restoring column F to the workbook's original (default/general) styling.
$ws.Columns.Item(6).ClearFormats()
